$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.690000000000001
$ws.Range("C4").Value = -11.684
$ws.Range("B6").Value = 6.161
$ws.Range("D6").Value = -7.237
$ws.Range("B7").Value = 5.717000000000001
$ws.Range("D7").Value = -7.589000000000001
$ws.Range("B8").Value = 5.988999999999999
$ws.Range("C8").Value = -11.682
$ws.Range("D8").Value = -7.741
$ws.Range("C9").Value = -10.827
$ws.Range("D10").Value = -7.903
$ws.Range("C12").Value = -11.053
$ws.Range("D13").Value = -7.947000000000001
$ws.Range("D14").Value = -8.035
$ws.Range("B16").Value = 5.239
$ws.Range("D16").Value = -7.965000000000001
$ws.Range("C17").Value = -13.011
$ws.Range("C18").Value = -12.646
$ws.Range("C19").Value = -11.949
$ws.Range("B20").Value = 8.607000000000001
$ws.Range("C20").Value = -12.23
$ws.Range("B21").Value = 8.580000000000002
$ws.Range("C26").Value = -12.635
$ws.Range("B28").Value = 5.393000000000001
$ws.Range("B29").Value = 5.784000000000001
$ws.Range("B30").Value = 5.217000000000001
$ws.Range("D30").Value = -7.081
$ws.Range("C31").Value = -12.423
$ws.Range("B32").Value = 6.986
$ws.Range("D37").Value = -8.065
$ws.Range("C39").Value = -12.059
$ws.Range("B40").Value = 9.023999999999999
$ws.Range("C40").Value = -12.09
$ws.Range("D40").Value = -8.430000000000001
$ws.Range("C41").Value = -12.02
$ws.Range("C42").Value = -12.364
$ws.Range("C43").Value = -12.182
$ws.Range("D44").Value = -7.742999999999999
$ws.Range("B46").Value = 5.299
$ws.Range("C47").Value = -12.476
$ws.Range("C48").Value = -11.844
$ws.Range("B51").Value = 5.398999999999999
$ws.Range("B52").Value = 5.628
$ws.Range("C54").Value = -12.867
$ws.Range("B57").Value = 4.880000000000001
$ws.Range("B59").Value = 5.086
$ws.Range("B62").Value = 5.146
$ws.Range("C62").Value = -13.345
$ws.Range("C63").Value = -11.174
$ws.Range("C64").Value = -11.151
$ws.Range("B66").Value = 5.412000000000001
$ws.Range("D70").Value = -7.567
$ws.Range("B73").Value = 6.102000000000001
$ws.Range("B74").Value = 9.186999999999999
$ws.Range("C76").Value = -12.17
$ws.Range("B77").Value = 6.506
$ws.Range("C81").Value = -12.378
$ws.Range("C84").Value = -13.466
$ws.Range("C89").Value = -13.233
$ws.Range("D89").Value = -8.32
$ws.Range("D91").Value = -7.486
$ws.Range("B92").Value = 4.679
$ws.Range("D93").Value = -7.331000000000001
$ws.Range("C94").Value = -11.931
$ws.Range("D98").Value = -7.154999999999999
$ws.Range("B100").Value = 6.121
